# Update the cached "datetimeFigureOut" date field text from 9/28/2021 to
# 9/30/2021 everywhere it appears: the slide master and every slide layout
# (each layout keeps its own cached copy of the date placeholder text).

$p = $ppt.ActivePresentation
$newDate = "9/30/2021"
$ppPlaceholderDate = 16

function Update-DateShape {
    param($shapes)
    for ($k = 1; $k -le $shapes.Count; $k++) {
        $shp = $shapes.Item($k)
        if ($shp.Type -eq 14) {
            if ($shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                if ($shp.HasTextFrame -eq -1) {
                    $shp.TextFrame.TextRange.Text = $newDate
                }
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateShape $master.Shapes

for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    Update-DateShape $layout.Shapes
}
